# Introduction.pptx - "Added missing names to js basics presentation."
#
# 1) Slide 15 ("contact"): append two new name paragraphs after the
#    existing "Spuzic Nenad" entry: "Janko Sokolovic" and "Ognjen Kurtic".
# 2) Slide 9 ("JS functions"): the description paragraph was re-typed by
#    the author, which collapsed its three separate runs into a single
#    run with identical combined text.

$p = $ppt.ActivePresentation

# --- 1. Contact slide: add the two missing names ------------------------
$contactSlide = $p.Slides.Item(15)
$contactShape = $contactSlide.Shapes.Item(2)
$contactRange = $contactShape.TextFrame.TextRange

$contactRange.InsertAfter("`rJanko Sokolovic`rOgnjen Kurtic")

# --- 2. JS functions slide: merge the split description runs -----------
$jsSlide = $p.Slides.Item(9)
$jsShape = $jsSlide.Shapes.Item(2)
$jsRange = $jsShape.TextFrame.TextRange
$jsPara = $jsRange.Paragraphs(1)

# Force a real text change so the run-splitting collapses into one run
# (re-assigning the exact same string is treated as a no-op).
$jsPara.Text = "__tmp_merge_placeholder__"
$jsPara = $jsRange.Paragraphs(1)
$jsPara.Text = "JS functions presentation covers how functions can be declared, gives an overview of closures, scopes and higher order functions."
